$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "51.148.64"
$ws.Range("E2").Value = "  -0.74%  "
$ws.Range("D3").Value = "2.957.69"
$ws.Range("E3").Value = "  +0.71%  "
$ws.Range("E4").Value = "  +0.02%  "
$ws.Range("D5").Value = "'380.42"
$ws.Range("E5").Value = "  +0.73%  "
$ws.Range("D6").Value = "'103.00"
$ws.Range("E6").Value = "  -0.72%  "
$ws.Range("D7").Value = "'0.539"
$ws.Range("E7").Value = "  -0.53%  "
$ws.Range("E8").Value = "  -0.02%  "
$ws.Range("D9").Value = "'0.590"
$ws.Range("E9").Value = "  +0.53%  "
$ws.Range("E10").Value = "  -1.02%  "
$ws.Range("E11").Value = "  -0.19%  "
$ws.Range("E12").Value = "  +0.18%  "
$ws.Range("D13").Value = "3.421.93"
$ws.Range("E13").Value = "  +0.60%  "
$ws.Range("D14").Value = "'18.09"
$ws.Range("E14").Value = "  -1.72%  "
$ws.Range("E15").Value = "  +0.53%  "
$ws.Range("D16").Value = "2.944.95"
$ws.Range("E16").Value = "  -0.02%  "
$ws.Range("D17").Value = "'0.987"
$ws.Range("E17").Value = "  +4.68%  "
$ws.Range("D18").Value = "51.109.97"
$ws.Range("E18").Value = "  -0.79%  "
$ws.Range("E19").Value = "  -5.75%  "
$ws.Range("D20").Value = "'7.08"
$ws.Range("E20").Value = "  -3.39%  "
$ws.Range("D21").Value = "'12.55"
$ws.Range("E21").Value = "  -3.53%  "
$ws.Range("D22").Value = "0.0₃0954"
$ws.Range("E22").Value = "  +0.43%  "
$ws.Range("D23").Value = "'68.51"
$ws.Range("E23").Value = "  +0.20%  "
$ws.Range("D24").Value = "'261.87"
$ws.Range("E24").Value = "  -0.02%  "
$ws.Range("D25").Value = "'2.87"
$ws.Range("E25").Value = "  +1.64%  "
$ws.Range("D26").Value = "'8.40"
$ws.Range("E26").Value = "  +14.76%  "
$ws.Range("D27").Value = "'7.62"
$ws.Range("E27").Value = "  +5.79%  "
$ws.Range("E28").Value = "  +0.67%  "
$ws.Range("E29").Value = "  +9.75%  "
$ws.Range("D30").Value = "'4.09"
$ws.Range("E30").Value = "  -1.03%  "
$ws.Range("E31").Value = "  -0.04%  "
$ws.Range("D32").Value = "'25.70"
$ws.Range("E32").Value = "  -0.53%  "
$ws.Range("D33").Value = "'9.82"
$ws.Range("E33").Value = "  -0.19%  "
$ws.Range("D34").Value = "'0.0456"
$ws.Range("E34").Value = "  +6.37%  "
$ws.Range("D35").Value = "'33.90"
$ws.Range("E35").Value = "  -0.60%  "
$ws.Range("B36").Value = "Toncoin"
$ws.Range("C36").Value = "https://coinranking.com/coin/67YlI0K1b+toncoin-ton"
$ws.Range("D36").Value = "'2.05"
$ws.Range("E36").Value = "  -1.96%  "
$ws.Range("B37").Value = "OKB"
$ws.Range("C37").Value = "https://coinranking.com/coin/PDKcptVnzJTmN+okb-okb"
$ws.Range("D37").Value = "'50.32"
$ws.Range("E37").Value = "  -2.80%  "
$ws.Range("E38").Value = "  -0.05%  "
$ws.Range("D39").Value = "'2.98"
$ws.Range("E39").Value = "  -1.58%  "
$ws.Range("E40").Value = "  -1.42%  "
$ws.Range("E41").Value = "  -1.40%  "
$ws.Range("E42").Value = "  +0.42%  "
$ws.Range("E43").Value = "  -2.40%  "
$ws.Range("D44").Value = "'121.66"
$ws.Range("E44").Value = "  -2.43%  "
$ws.Range("D45").Value = "'21.10"
$ws.Range("E45").Value = "  -3.36%  "
$ws.Range("D46").Value = "'2.07"
$ws.Range("E46").Value = "  +0.27%  "
$ws.Range("E47").Value = "  -0.67%  "
$ws.Range("E48").Value = "  +2.38%  "
$ws.Range("D49").Value = "2.007.87"
$ws.Range("E49").Value = "  -0.75%  "
$ws.Range("E50").Value = "  +1.81%  "
$ws.Range("D51").Value = "'0.0336"
$ws.Range("E51").Value = "  +4.17%  "
